$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2063.0889
$ws.Range("I15").Value = 2063.0889
$ws.Range("K15").Value = 6189.2667
$ws.Range("M15").Value = -6020.2667

# Row 132
$ws.Range("H132").Value = 2086609.4
$ws.Range("I132").Value = 3139.558
$ws.Range("J132").Value = 20004450
$ws.Range("K132").Value = 9418.673999999999
$ws.Range("L132").Value = 60013350
$ws.Range("M132").Value = -6888.673999999999
$ws.Range("N132").Value = -60018410

# Row 137
$ws.Range("H137").Value = 9856.538
$ws.Range("I137").Value = 17549.385
$ws.Range("K137").Value = 52648.155
$ws.Range("M137").Value = -50098.155

# Row 138
$ws.Range("H138").Value = 229641.64
$ws.Range("J138").Value = 4037.5312
$ws.Range("L138").Value = 12112.5936
$ws.Range("N138").Value = -22392.5936

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9787.725
$ws.Range("I32").Value = 9782.281999999999
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 9782.281999999999
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -9495.281999999999
$ws.Range("N32").Value = -10574

# Row 37
$ws.Range("H37").Value = 44666.668
$ws.Range("I37").Value = 34000
$ws.Range("K37").Value = 34000
$ws.Range("M37").Value = -33727

# Row 74
$ws.Range("H74").Value = 6669.4546
$ws.Range("I74").Value = 7732.9375
$ws.Range("J74").Value = 3833.5
$ws.Range("K74").Value = 7732.9375
$ws.Range("L74").Value = 3833.5
$ws.Range("M74").Value = -6858.9375
$ws.Range("N74").Value = -5581.5

# Row 77
$ws.Range("H77").Value = 6669.4546
$ws.Range("I77").Value = 7732.9375
$ws.Range("J77").Value = 3833.5
$ws.Range("K77").Value = 38664.6875
$ws.Range("L77").Value = 19167.5
$ws.Range("M77").Value = -34296.6875
$ws.Range("N77").Value = -27903.5

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 222.14285
$ws.Range("I22").Value = 233.33333
$ws.Range("K22").Value = 233.33333
$ws.Range("M22").Value = -60.33332999999999

# Row 86
$ws.Range("H86").Value = 4672.3335
$ws.Range("I86").Value = 7283.0713
$ws.Range("J86").Value = 1860.7693
$ws.Range("K86").Value = 7283.0713
$ws.Range("L86").Value = 1860.7693
$ws.Range("M86").Value = -6160.0713
$ws.Range("N86").Value = -4106.7693

# Row 89
$ws.Range("H89").Value = 4672.3335
$ws.Range("I89").Value = 7283.0713
$ws.Range("J89").Value = 1860.7693
$ws.Range("K89").Value = 36415.35649999999
$ws.Range("L89").Value = 9303.8465
$ws.Range("M89").Value = -30799.35649999999
$ws.Range("N89").Value = -20535.8465

# Row 105
$ws.Range("H105").Value = 131499
$ws.Range("I105").Value = 252748
$ws.Range("K105").Value = 252748
$ws.Range("M105").Value = -251001

# Row 134
$ws.Range("H134").Value = 6320.654
$ws.Range("I134").Value = 6993.9546
$ws.Range("J134").Value = 2617.5
$ws.Range("K134").Value = 20981.8638
$ws.Range("L134").Value = 7852.5
$ws.Range("M134").Value = -18446.8638
$ws.Range("N134").Value = -12922.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 11181.066
$ws.Range("I31").Value = 12645.25
$ws.Range("J31").Value = 5324.3335
$ws.Range("K31").Value = 12645.25
$ws.Range("L31").Value = 5324.3335
$ws.Range("M31").Value = -12350.25
$ws.Range("N31").Value = -5914.3335

# Row 34
$ws.Range("H34").Value = 11181.066
$ws.Range("I34").Value = 12645.25
$ws.Range("J34").Value = 5324.3335
$ws.Range("K34").Value = 12645.25
$ws.Range("L34").Value = 5324.3335
$ws.Range("M34").Value = -12443.25
$ws.Range("N34").Value = -5728.3335

# Row 62
$ws.Range("H62").Value = 6646.8
$ws.Range("I62").Value = 6369.8335
$ws.Range("J62").Value = 7062.25
$ws.Range("K62").Value = 6369.8335
$ws.Range("L62").Value = 7062.25
$ws.Range("M62").Value = -5745.8335
$ws.Range("N62").Value = -8310.25

# Row 65
$ws.Range("H65").Value = 6646.8
$ws.Range("I65").Value = 6369.8335
$ws.Range("J65").Value = 7062.25
$ws.Range("K65").Value = 31849.1675
$ws.Range("L65").Value = 35311.25
$ws.Range("M65").Value = -28729.1675
$ws.Range("N65").Value = -41551.25

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 105
$ws.Range("H105").Value = 423789.4
$ws.Range("I105").Value = 528486.75
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 528486.75
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -526739.75
$ws.Range("N105").Value = -8494

# Row 132
$ws.Range("H132").Value = 2441.3845
$ws.Range("I132").Value = 2567.0908
$ws.Range("K132").Value = 7701.2724
$ws.Range("M132").Value = -5171.2724

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 343.35294
$ws.Range("I33").Value = 388.8
$ws.Range("J33").Value = 278.42856
$ws.Range("K33").Value = 2332.8
$ws.Range("L33").Value = 1670.57136
$ws.Range("M33").Value = -2049.8
$ws.Range("N33").Value = -2236.57136

# Row 113
$ws.Range("H113").Value = 14147.125
$ws.Range("J113").Value = 14147.125
$ws.Range("L113").Value = 42441.375
$ws.Range("N113").Value = -46781.375

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9619.091
$ws.Range("I80").Value = 12608.4
$ws.Range("J80").Value = 7128
$ws.Range("K80").Value = 12608.4
$ws.Range("L80").Value = 7128
$ws.Range("M80").Value = -11610.4
$ws.Range("N80").Value = -9124

# Row 83
$ws.Range("H83").Value = 9619.091
$ws.Range("I83").Value = 12608.4
$ws.Range("J83").Value = 7128
$ws.Range("K83").Value = 63042
$ws.Range("L83").Value = 35640
$ws.Range("M83").Value = -58050
$ws.Range("N83").Value = -45624

# Row 102
$ws.Range("H102").Value = 7677.619
$ws.Range("I102").Value = 9222.875
$ws.Range("K102").Value = 9222.875
$ws.Range("M102").Value = -7600.875

# Row 122
$ws.Range("H122").Value = 11557
$ws.Range("I122").Value = 7723.4165
$ws.Range("K122").Value = 23170.2495
$ws.Range("M122").Value = -20720.2495

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 36285.43
$ws.Range("I7").Value = 47599.6
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 47599.6
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -47487.6
$ws.Range("N7").Value = -8224

# Row 40
$ws.Range("H40").Value = 17620.527
$ws.Range("I40").Value = 19834.5
$ws.Range("K40").Value = 19834.5
$ws.Range("M40").Value = -19698.5

# Row 68
$ws.Range("H68").Value = 3776.5
$ws.Range("I68").Value = 2075.3333
$ws.Range("J68").Value = 8880
$ws.Range("K68").Value = 2075.3333
$ws.Range("L68").Value = 8880
$ws.Range("M68").Value = -1326.3333
$ws.Range("N68").Value = -10378

# Row 71
$ws.Range("H71").Value = 3776.5
$ws.Range("I71").Value = 2075.3333
$ws.Range("J71").Value = 8880
$ws.Range("K71").Value = 10376.6665
$ws.Range("L71").Value = 44400
$ws.Range("M71").Value = -6632.666499999999
$ws.Range("N71").Value = -51888

# Row 126
$ws.Range("H126").Value = 36285.43
$ws.Range("I126").Value = 47599.6
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 142798.8
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -140328.8
$ws.Range("N126").Value = -28940

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 17628.074
$ws.Range("I132").Value = 30625.072
$ws.Range("J132").Value = 3631.3076
$ws.Range("K132").Value = 91875.216
$ws.Range("L132").Value = 10893.9228
$ws.Range("M132").Value = -89345.216
$ws.Range("N132").Value = -15953.9228
